$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1162.4286
$ws.Range("I18").Value = 1094.3704
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 1094.3704
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -810.3704
$ws.Range("N18").Value = -3568

$ws.Range("H53").Value = 56068.945
$ws.Range("I53").Value = 166896.83
$ws.Range("J53").Value = 655
$ws.Range("K53").Value = 166896.83
$ws.Range("L53").Value = 655
$ws.Range("M53").Value = -166259.83
$ws.Range("N53").Value = -1929

$ws.Range("H62").Value = 4255.8335
$ws.Range("I62").Value = 4255.8335
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4255.8335
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3631.8335
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 4255.8335
$ws.Range("I65").Value = 4255.8335
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 21279.1675
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -18159.1675
$ws.Range("N65").ClearContents()

$ws.Range("H86").Value = 500751.5
$ws.Range("I86").Value = 1000003
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1000003
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -998880
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 500751.5
$ws.Range("I89").Value = 1000003
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 5000015
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -4994399
$ws.Range("N89").Value = -18732

$ws.Range("H129").Value = 4330.033
$ws.Range("I129").Value = 723.6667
$ws.Range("J129").Value = 4730.7407
$ws.Range("K129").Value = 2171.0001
$ws.Range("L129").Value = 14192.2221
$ws.Range("M129").Value = 2828.9999
$ws.Range("N129").Value = -24192.2221

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 833.8461
$ws.Range("I122").Value = 576.36365
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 1729.09095
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = 720.90905
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 196
$ws.Range("I7").Value = 49
$ws.Range("J7").Value = 490
$ws.Range("K7").Value = 147
$ws.Range("L7").Value = 1470
$ws.Range("M7").Value = -35
$ws.Range("N7").Value = -1694

$ws.Range("H64").Value = 1409.3334
$ws.Range("J64").Value = 1669
$ws.Range("L64").Value = 5007
$ws.Range("N64").Value = -5547

$ws.Range("H67").Value = 1409.3334
$ws.Range("J67").Value = 1669
$ws.Range("L67").Value = 5007
$ws.Range("N67").Value = -6879

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H129").Value = 1379.1177
$ws.Range("J129").Value = 1543.9286
$ws.Range("L129").Value = 4631.7858
$ws.Range("N129").Value = -14631.7858

$ws.Range("H137").Value = 26330916
$ws.Range("I137").Value = 1196.9231
$ws.Range("J137").Value = 83378640
$ws.Range("K137").Value = 3590.7693
$ws.Range("L137").Value = 250135920
$ws.Range("M137").Value = 1509.2307
$ws.Range("N137").Value = -250146120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3501.3784
$ws.Range("I80").Value = 2791.4644
$ws.Range("J80").Value = 5710
$ws.Range("K80").Value = 2791.4644
$ws.Range("L80").Value = 5710
$ws.Range("M80").Value = -1793.4644
$ws.Range("N80").Value = -7706

$ws.Range("H83").Value = 3501.3784
$ws.Range("I83").Value = 2791.4644
$ws.Range("J83").Value = 5710
$ws.Range("K83").Value = 13957.322
$ws.Range("L83").Value = 28550
$ws.Range("M83").Value = -8965.322
$ws.Range("N83").Value = -38534

$ws.Range("H102").Value = 1509.1818
$ws.Range("I102").Value = 1800
$ws.Range("J102").Value = 1400.125
$ws.Range("K102").Value = 1800
$ws.Range("L102").Value = 1400.125
$ws.Range("M102").Value = -178
$ws.Range("N102").Value = -4644.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 655.0625
$ws.Range("I22").Value = 716.6667
$ws.Range("J22").Value = 618.1
$ws.Range("K22").Value = 716.6667
$ws.Range("L22").Value = 618.1
$ws.Range("M22").Value = -421.6667
$ws.Range("N22").Value = -1208.1

$ws.Range("H27").Value = 655.0625
$ws.Range("I27").Value = 716.6667
$ws.Range("J27").Value = 618.1
$ws.Range("K27").Value = 716.6667
$ws.Range("L27").Value = 618.1
$ws.Range("M27").Value = -609.6667
$ws.Range("N27").Value = -832.1

$ws.Range("H55").Value = 656.25
$ws.Range("I55").Value = 490
$ws.Range("J55").Value = 933.3333
$ws.Range("K55").Value = 490
$ws.Range("L55").Value = 933.3333
$ws.Range("M55").Value = -317
$ws.Range("N55").Value = -1279.3333

$ws.Range("H122").Value = 3759.5386
$ws.Range("I122").Value = 4097.1113
$ws.Range("K122").Value = 12291.3339
$ws.Range("M122").Value = -9841.333899999998

$ws.Range("H132").Value = 2882.1853
$ws.Range("I132").Value = 2403.2354
$ws.Range("K132").Value = 7209.706200000001
$ws.Range("M132").Value = -4679.706200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1014.9286
$ws.Range("I126").Value = 1137.1818
$ws.Range("J126").Value = 566.6667
$ws.Range("K126").Value = 3411.5454
$ws.Range("L126").Value = 1700.0001
$ws.Range("M126").Value = -941.5454
$ws.Range("N126").Value = -6640.0001

$ws.Range("H132").Value = 1204.3
$ws.Range("I132").Value = 760.94116
$ws.Range("J132").Value = 1784.0769
$ws.Range("K132").Value = 2282.82348
$ws.Range("L132").Value = 5352.2307
$ws.Range("M132").Value = 247.17652
$ws.Range("N132").Value = -10412.2307
